# Updates cryptos list values (prices / 1h volume %) and reorders two
# coin-pair rows, per the "Updated cryptos list ... with GitHub Actions" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '45.266.54'
$ws.Range("E2").Value = '  +4.13%  '
# Row 3
$ws.Range("D3").Value = '2.369.85'
$ws.Range("E3").Value = '  +1.92%  '
# Row 4
$ws.Range("E4").Value = '  +0.07%  '
# Row 5
$ws.Range("D5").Value = '''109.00'
$ws.Range("E5").Value = '  +0.56%  '
# Row 6
$ws.Range("D6").Value = '''310.56'
$ws.Range("E6").Value = '  -0.12%  '
# Row 7
$ws.Range("E7").Value = '  +0.11%  '
# Row 8
$ws.Range("E8").Value = '  -0.16%  '
# Row 9
$ws.Range("E9").Value = '  +0.58%  '
# Row 10
$ws.Range("D10").Value = '''41.22'
$ws.Range("E10").Value = '  +0.91%  '
# Row 11
$ws.Range("E11").Value = '  +0.29%  '
# Row 12
$ws.Range("D12").Value = '''8.49'
$ws.Range("E12").Value = '  -1.29%  '
# Row 13
$ws.Range("D13").Value = '''0.109'
$ws.Range("E13").Value = '  +1.45%  '
# Row 14
$ws.Range("E14").Value = '  -2.01%  '
# Row 15
$ws.Range("D15").Value = '2.731.53'
$ws.Range("E15").Value = '  +1.95%  '
# Row 16
$ws.Range("E16").Value = '  -1.07%  '
# Row 17
$ws.Range("D17").Value = '2.371.72'
$ws.Range("E17").Value = '  +1.80%  '
# Row 18
$ws.Range("D18").Value = '45.275.88'
$ws.Range("E18").Value = '  +4.42%  '
# Row 19
$ws.Range("D19").Value = '''14.98'
$ws.Range("E19").Value = '  +13.53%  '
# Row 20
$ws.Range("D20").Value = '''7.29'
$ws.Range("E20").Value = '  -3.47%  '
# Row 21
$ws.Range("E21").Value = '  -0.38%  '
# Row 22
$ws.Range("D22").Value = '''73.32'
$ws.Range("E22").Value = '  -1.07%  '
# Row 23
$ws.Range("E23").Value = '  +0.08%  '
# Row 24
$ws.Range("D24").Value = '''260.06'
$ws.Range("E24").Value = '  -3.02%  '
# Row 25
$ws.Range("E25").Value = '  +2.03%  '
# Row 26
$ws.Range("E26").Value = '  -0.29%  '
# Row 27
$ws.Range("D27").Value = '''11.15'
$ws.Range("E27").Value = '  +0.15%  '
# Row 28
$ws.Range("D28").Value = '''7.27'
$ws.Range("E28").Value = '  -4.81%  '
# Row 29
$ws.Range("E29").Value = '  +2.25%  '
# Row 30
$ws.Range("D30").Value = '''0.0969'
$ws.Range("E30").Value = '  +9.36%  '
# Row 31
$ws.Range("D31").Value = '''22.43'
$ws.Range("E31").Value = '  -0.79%  '
# Row 32
$ws.Range("D32").Value = '''37.77'
$ws.Range("E32").Value = '  -2.33%  '
# Row 33
$ws.Range("D33").Value = '''169.38'
$ws.Range("E33").Value = '  +1.41%  '
# Row 34
$ws.Range("E34").Value = '  +4.39%  '
# Row 35
$ws.Range("E35").Value = '  -0.34%  '
# Row 36
$ws.Range("E36").Value = '  +4.02%  '
# Row 37
$ws.Range("D37").Value = '''4.74'
$ws.Range("E37").Value = '  +0.06%  '
# Row 38
$ws.Range("B38").Value = 'LidoDAOToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D38").Value = '''2.96'
$ws.Range("E38").Value = '  +4.49%  '
# Row 39
$ws.Range("B39").Value = 'NEARProtocol'
$ws.Range("C39").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D39").Value = '''3.94'
$ws.Range("E39").Value = '  +3.79%  '
# Row 40
$ws.Range("D40").Value = '''0.0354'
$ws.Range("E40").Value = '  -2.46%  '
# Row 41
$ws.Range("E41").Value = '  +3.65%  '
# Row 42
$ws.Range("D42").Value = '''99.68'
$ws.Range("E42").Value = '  -4.66%  '
# Row 43
$ws.Range("D43").Value = '''0.230'
$ws.Range("E43").Value = '  -2.91%  '
# Row 44
$ws.Range("D44").Value = '''69.47'
$ws.Range("E44").Value = '  -3.13%  '
# Row 45
$ws.Range("D45").Value = '''12.94'
$ws.Range("E45").Value = '  -2.56%  '
# Row 46
$ws.Range("E46").Value = '  -0.32%  '
# Row 47
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '1.759.71'
$ws.Range("E47").Value = '  +6.04%  '
# Row 48
$ws.Range("B48").Value = 'ordi'
$ws.Range("C48").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D48").Value = '''81.54'
$ws.Range("E48").Value = '  +9.10%  '
# Row 49
$ws.Range("D49").Value = '''5.56'
# Row 50
$ws.Range("D50").Value = '''111.96'
$ws.Range("E50").Value = '  -1.71%  '
# Row 51
$ws.Range("D51").Value = '''9.18'
$ws.Range("E51").Value = '  +2.75%  '
